$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.054.68'
$ws.Range('E2').Value = '  +2.49%  '

$ws.Range('D3').Value = '3.111.85'
$ws.Range('E3').Value = '  +2.50%  '

$ws.Range('E4').Value = '  -0.09%  '

$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '543.32'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  +0.95%  '

$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '139.53'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  +5.28%  '

$ws.Range('E7').Value = '  -0.13%  '

$ws.Range('D8').Value = '3.102.40'
$ws.Range('E8').Value = '  +2.58%  '

$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.499'
$cell.Style = "Normal"
$ws.Range('E9').Value = '  +3.45%  '

$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '0.158'
$cell.Style = "Normal"
$ws.Range('E10').Value = '  +3.42%  '

$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '6.52'
$cell.Style = "Normal"
$ws.Range('E11').Value = '  +2.40%  '

$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '0.462'
$cell.Style = "Normal"
$ws.Range('E12').Value = '  +2.21%  '

$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '0.0000229'
$cell.Style = "Normal"
$ws.Range('E13').Value = '  +8.34%  '

$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '35.10'
$cell.Style = "Normal"
$ws.Range('E14').Value = '  +2.71%  '

$ws.Range('D15').Value = '3.602.77'
$ws.Range('E15').Value = '  +2.04%  '

$ws.Range('D16').Value = '64.028.75'
$ws.Range('E16').Value = '  +2.29%  '

$ws.Range('E17').Value = '  +2.21%  '

$ws.Range('D18').Value = '3.099.86'
$ws.Range('E18').Value = '  +2.16%  '

$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '6.74'
$cell.Style = "Normal"
$ws.Range('E19').Value = '  +3.41%  '

$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '486.79'
$cell.Style = "Normal"
$ws.Range('E20').Value = '  +2.24%  '

$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '13.52'
$cell.Style = "Normal"
$ws.Range('E21').Value = '  +2.65%  '

$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '0.706'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  +2.94%  '

$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '7.18'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  +3.44%  '

$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '79.60'
$cell.Style = "Normal"
$ws.Range('E24').Value = '  +3.95%  '

$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '12.34'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  +3.08%  '

$ws.Range('E26').Value = '  +0.26%  '

$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '2.74'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  +2.91%  '

$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '8.20'
$cell.Style = "Normal"
$ws.Range('E28').Value = '  +0.77%  '

$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '0.997'
$cell.Style = "Normal"
$ws.Range('E29').Value = '  -0.27%  '

$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '26.50'
$cell.Style = "Normal"
$ws.Range('E30').Value = '  +2.63%  '

$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '1.92'
$cell.Style = "Normal"
$ws.Range('E31').Value = '  +1.47%  '

$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '1.15'
$cell.Style = "Normal"
$ws.Range('E32').Value = '  +3.35%  '

$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '57.64'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  -3.53%  '

$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '2.38'
$cell.Style = "Normal"
$ws.Range('E34').Value = '  -2.45%  '

$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '502.81'
$cell.Style = "Normal"
$ws.Range('E35').Value = '  -0.76%  '

$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '5.38'
$cell.Style = "Normal"
$ws.Range('E36').Value = '  +6.73%  '

$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '6.06'
$cell.Style = "Normal"
$ws.Range('E37').Value = '  +4.26%  '

$ws.Range('D38').Value = '3.273.46'
$ws.Range('E38').Value = '  +7.72%  '

$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '0.0408'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  +3.78%  '

$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '0.0802'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  +3.42%  '

$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.119'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  +3.42%  '

$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '2.72'
$cell.Style = "Normal"
$ws.Range('E42').Value = '  +6.50%  '

$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '8.17'
$cell.Style = "Normal"
$ws.Range('E43').Value = '  +3.06%  '

$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '0.258'
$cell.Style = "Normal"
$ws.Range('E44').Value = '  +4.26%  '

$ws.Range('E45').Value = '  +0.05%  '

$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '2.07'
$cell.Style = "Normal"
$ws.Range('E46').Value = '  +3.93%  '

$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '123.84'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  +3.58%  '

$ws.Range('D48').Value = '0.0₃0538'
$ws.Range('E48').Value = '  +10.95%  '

$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '24.88'
$cell.Style = "Normal"
$ws.Range('E49').Value = '  +4.62%  '

$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '0.109'
$cell.Style = "Normal"
$ws.Range('E50').Value = '  +3.50%  '

$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '2.41'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  +3.21%  '
